$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window position (workbook.xml bookViews) ---
$wb.Windows.Item(1).Left = -38000
$wb.Windows.Item(1).Top = 1360

# --- Header row (row 2) label changes ---
$ws.Range("B2").Value2 = "in:wby.ship:refrigeration:set"
$ws.Range("C2").Value2 = "in:ProSpec quantity"
$ws.Range("D2").Value2 = "in:ProSpec total"
$ws.Range("E2").Value2 = "out:ProSpec Shipping Method"
$ws.Range("G2").Value2 = "out:ProSpec Handling"
$ws.Range("H2").Value2 = "out:FedExZone 1"
$ws.Range("I2").Value2 = "out:FedExZone 2"
$ws.Range("J2").Value2 = "out:FedExZone 3"
$ws.Range("K2").Value2 = "out:FedExZone 4"
$ws.Range("L2").Value2 = "out:FedExZone 5"

# --- Rule columns A/B for rows 3-8 ---
$ws.Range("A3").Value2 = "`$(ProSpec) <=  `$in"
$ws.Range("B3").Value2 = "`$in >= `$(ice-packs)"

$ws.Range("A4").Value2 = "`$(ProSpec) <=  `$in"
$ws.Range("B4").Value2 = "`$in >= `$(ice-packs)"

$ws.Range("A5").Value2 = "`$(ProSpec) <=  `$in"
$ws.Range("B5").Value2 = "`$in >= `$(ice-packs)"

$ws.Range("A6").Value2 = "`$(ProSpec) <=  `$in"
$ws.Range("B6").Value2 = "`$in >= `$(ice-packs)"

$ws.Range("A7").Value2 = "`$(ProSpec) <=  `$in"

$ws.Range("A8").Value2 = "`$(ProSpec) <=  `$in"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 24.6640625
$ws.Columns.Item(3).ColumnWidth = 20.1640625
$ws.Columns.Item(6).ColumnWidth = 19.5
$ws.Columns.Item(7).ColumnWidth = 23
$ws.Columns.Item(8).ColumnWidth = 16.1640625
$ws.Columns.Item(9).ColumnWidth = 14.1640625
$ws.Columns.Item(10).ColumnWidth = 14.1640625
$ws.Columns.Item(11).ColumnWidth = 15.6640625
$ws.Columns.Item(12).ColumnWidth = 14.1640625

# --- Selection ---
$ws.Range("K2").Select()
